$d = $word.ActiveDocument

# Replace the whole text of a paragraph while preserving its paragraph /
# cell-mark, and without disturbing the run's xml:space="preserve"
# attribute the way Find.Execute's "Replace With" / a Find-range .Text
# assignment would (that path always re-serializes <w:t> without
# xml:space, even when the original had it).
function Set-ParaText($findText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $r = $p.Range
        $t = $r.Text
        # Paragraphs inside table cells end with a cell-mark (chr 7) in
        # addition to / instead of the paragraph mark (chr 13); strip
        # both so we can match on the visible text only.
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $findText) {
            $r.End = $r.Start + $trimmed.Length
            $r.Text = $newText
            return $true
        }
    }
    return $false
}

Set-ParaText "N. d'étudiant·e·s" "የ d'étudiant·e·s" | Out-Null
Set-ParaText "Date" "ቀን" | Out-Null
Set-ParaText "00:00 - 00:26" "Video Title" | Out-Null
Set-ParaText "Conclusion" "ማጠቃለያ" | Out-Null
